# "Generate Report for Handoff": refresh Priority + Latest Handoff Datetime
# for the four in-flight files on the zh-cn and de-de handback sheets.

$wb = $excel.ActiveWorkbook

$zh = $wb.Worksheets.Item("zh-cn")
$de = $wb.Worksheets.Item("de-de")

# zh-cn: rows 4-7, column E (Priority) low -> ht, column H (Latest Handoff Datetime) bumped
$zh.Range("E4").Value = "ht"
$zh.Range("H4").Value = "2016-09-04 06:34:25"

$zh.Range("E5").Value = "ht"
$zh.Range("H5").Value = "2016-09-04 06:34:25"

$zh.Range("E6").Value = "ht"
$zh.Range("H6").Value = "2016-09-04 06:34:25"

$zh.Range("E7").Value = "ht"
$zh.Range("H7").Value = "2016-09-04 06:34:25"

# de-de: rows 4-7, same Priority change, Latest Handoff Datetime bumped
$de.Range("E4").Value = "ht"
$de.Range("H4").Value = "2016-09-04 06:34:29"

$de.Range("E5").Value = "ht"
$de.Range("H5").Value = "2016-09-04 06:34:29"

$de.Range("E6").Value = "ht"
$de.Range("H6").Value = "2016-09-04 06:34:29"

$de.Range("E7").Value = "ht"
$de.Range("H7").Value = "2016-09-04 06:34:29"
